$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sorted data (descending by 2003 value), with Uzbek and Vietnamese removed
$data = @(
    @("English", 26.42567604536709),
    @("Chinese", 9.923786660264328),
    @("Spanish", 7.485259813021717),
    @("Japanese", 6.406387452471833),
    @("German", 5.374744172966636),
    @("Arabic", 4.892333088956317),
    @("Portuguese", 3.553313379271298),
    @("Russian", 3.407480749208801),
    @("French", 3.35187214991849),
    @("Italian", 3.07553577632755),
    @("Malay-Indonesian", 2.63997148984114),
    @("Korean", 1.661523644901343),
    @("Persian", 1.585374132860658),
    @("Dutch", 1.554924567781548),
    @("Turkish", 1.342994434327395),
    @("Thai", 0.9768009538793023),
    @("Polish", 0.8660633163225131),
    @("Urdu", 0.8081234702046004),
    @("Swedish", 0.5005149393157337),
    @("Bengali", 0.4492155627484779)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (previously rows 22 and 23)
$ws.Range("A22:B23").Delete()
